$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.275.55"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").Value = "1.872.83"
$ws.Range("E3").Value = "  +3.12%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.46%  "

# Row 5
$ws.Range("D5").Value = "'312.32"
$ws.Range("E5").Value = "  -0.12%  "

# Row 6
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").Value = "'0.5029"
$ws.Range("E7").Value = "  -2.40%  "

# Row 8
$ws.Range("D8").Value = "'0.3942"
$ws.Range("E8").Value = "  -0.56%  "

# Row 9
$ws.Range("D9").Value = "'0.09922"
$ws.Range("E9").Value = "  +26.64%  "

# Row 10
$ws.Range("D10").Value = "'1.125"
$ws.Range("E10").Value = "  +0.97%  "

# Row 11
$ws.Range("D11").Value = "'41.29"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("D12").Value = "'6.459"
$ws.Range("E12").Value = "  +1.33%  "

# Row 13
$ws.Range("D13").Value = "'20.93"
$ws.Range("E13").Value = "  +1.85%  "

# Row 14
$ws.Range("D14").Value = "1.876.97"
$ws.Range("E14").Value = "  +3.77%  "

# Row 15
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "  +0.45%  "

# Row 16
$ws.Range("D16").Value = "'7.389"
$ws.Range("E16").Value = "  +0.46%  "

# Row 17
$ws.Range("D17").Value = "'0.00001139"
$ws.Range("E17").Value = "  +5.39%  "

# Row 18
$ws.Range("D18").Value = "'93.54"
$ws.Range("E18").Value = "  +0.57%  "

# Row 19
$ws.Range("D19").Value = "'0.06666"
$ws.Range("E19").Value = "  +1.32%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'17.40"
$ws.Range("E20").Value = "  +0.21%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("D22").Value = "'6.102"
$ws.Range("E22").Value = "  +1.25%  "

# Row 23
$ws.Range("D23").Value = "28.323.41"
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").Value = "'11.33"
$ws.Range("E24").Value = "  +1.43%  "

# Row 25
$ws.Range("D25").Value = "'2.264"
$ws.Range("E25").Value = "  +2.02%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.526"
$ws.Range("E26").Value = "  +2.38%  "

# Row 27
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.092.38"
$ws.Range("E27").Value = "  +3.68%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.28"
$ws.Range("E28").Value = "  +3.18%  "

# Row 29
$ws.Range("D29").Value = "'157.81"
$ws.Range("E29").Value = "  -1.88%  "

# Row 30
$ws.Range("D30").Value = "'127.48"
$ws.Range("E30").Value = "  -0.51%  "

# Row 31
$ws.Range("D31").Value = "'0.1063"
$ws.Range("E31").Value = "  -3.28%  "

# Row 32
$ws.Range("D32").Value = "'1.056"
$ws.Range("E32").Value = "  -1.05%  "

# Row 33
$ws.Range("D33").Value = "'5.637"
$ws.Range("E33").Value = "  +0.88%  "

# Row 34
$ws.Range("D34").Value = "'3.607"
$ws.Range("E34").Value = "  -1.27%  "

# Row 35
$ws.Range("D35").Value = "'0.06813"
$ws.Range("E35").Value = "  -5.16%  "

# Row 36
$ws.Range("D36").Value = "'9.410"
$ws.Range("E36").Value = "  +2.23%  "

# Row 37
$ws.Range("D37").Value = "'0.02392"
$ws.Range("E37").Value = "  +1.49%  "

# Row 38
$ws.Range("D38").Value = "'0.2186"
$ws.Range("E38").Value = "  +0.07%  "

# Row 39
$ws.Range("D39").Value = "'5.015"
$ws.Range("E39").Value = "  -0.69%  "

# Row 40
$ws.Range("D40").Value = "'11.47"
$ws.Range("E40").Value = "  -1.15%  "

# Row 41
$ws.Range("D41").Value = "'0.6299"
$ws.Range("E41").Value = "  +1.50%  "

# Row 42
$ws.Range("D42").Value = "'1.176"
$ws.Range("E42").Value = "  +1.12%  "

# Row 43
$ws.Range("E43").Value = "  +0.21%  "

# Row 44
$ws.Range("D44").Value = "'13.45"
$ws.Range("E44").Value = "  +2.23%  "

# Row 45
$ws.Range("D45").Value = "'0.6016"
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("E46").Value = "  -1.61%  "

# Row 47
$ws.Range("D47").Value = "'3.667"
$ws.Range("E47").Value = "  -1.97%  "

# Row 48
$ws.Range("D48").Value = "'125.36"
$ws.Range("E48").Value = "  -0.11%  "

# Row 49
$ws.Range("D49").Value = "'1.994"
$ws.Range("E49").Value = "  +3.51%  "

# Row 50
$ws.Range("D50").Value = "'1.198"
$ws.Range("E50").Value = "  -1.15%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.122"
$ws.Range("E51").Value = "  +4.58%  "
